$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (GitHub Actions update).
# Price cells that look numeric are entered with a leading apostrophe so
# Excel stores them as text (matching the sheet's existing inlineStr
# formatting, e.g. trailing zeros like "6.560" or "242.07") instead of
# silently coercing them to a Number and dropping formatting.
$ws.Range('D2').Value = '29.190.67'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.859.06'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('D4').Value = "'0.9994"
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'242.07"
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('D7').Value = "'0.9996"
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = "'0.07806"
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('D9').Value = "'0.3107"
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').Value = "'23.92"
$ws.Range('E10').Value = '  -4.08%  '
$ws.Range('D11').Value = "'0.07800"
$ws.Range('E11').Value = '  -3.97%  '
$ws.Range('D12').Value = '1.867.76'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = "'92.65"
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'5.125"
$ws.Range('E14').Value = '  -2.11%  '
$ws.Range('D15').Value = "'0.6911"
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').Value = "'6.560"
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D17').Value = "'0.000008458"
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '29.219.35'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').Value = "'250.43"
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Value = '2.112.64'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('E21').Value = '  -3.23%  '
$ws.Range('D22').Value = "'0.9994"
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = "'7.604"
$ws.Range('E23').Value = '  -1.02%  '
$ws.Range('D24').Value = "'0.9997"
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = "'0.1533"
$ws.Range('E25').Value = '  -3.03%  '
$ws.Range('D26').Value = "'160.83"
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').Value = "'8.900"
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').Value = "'18.56"
$ws.Range('E28').Value = '  -2.23%  '
$ws.Range('D29').Value = "'1.569"
$ws.Range('E29').Value = '  +4.00%  '
$ws.Range('D30').Value = "'4.274"
$ws.Range('E30').Value = '  -3.34%  '
$ws.Range('D31').Value = "'4.255"
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').Value = "'1.213"
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('D33').Value = "'0.05233"
$ws.Range('E33').Value = '  -1.60%  '
$ws.Range('D34').Value = "'0.7583"
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = "'1.877"
$ws.Range('E35').Value = '  -3.70%  '
$ws.Range('D36').Value = "'1.177"
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = "'2.707"
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').Value = '1.222.02'
$ws.Range('E39').Value = '  -4.39%  '
$ws.Range('D40').Value = "'2.722"
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = "'0.8998"
$ws.Range('E41').Value = '  -0.98%  '
$ws.Range('D42').Value = "'110.48"
$ws.Range('E42').Value = '  -1.27%  '
$ws.Range('D43').Value = "'5.818"
$ws.Range('E43').Value = '  -9.21%  '
$ws.Range('D44').Value = "'0.9994"
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = "'67.44"
$ws.Range('E45').Value = '  -9.22%  '
$ws.Range('D46').Value = '2.008.56'
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('D47').Value = "'0.5185"
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = "'9.539"
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  -6.08%  '
$ws.Range('D50').Value = "'1.768"
$ws.Range('E50').Value = '  -2.18%  '
$ws.Range('D51').Value = "'7.040"
$ws.Range('E51').Value = '  -0.91%  '
